$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AD1:AF1").Font.Bold = $true
$ws.Range("AD1:AF1").HorizontalAlignment = -4108
$ws.Range("AD1:AF1").VerticalAlignment = -4160
$ws.Range("AD1:AF1").Borders.LineStyle = 1

# Data rows 2-50
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 30).Value = 65   # AD
    $ws.Cells.Item($r, 31).Value = 97   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
